$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 42606.881180555552

$ws.Range("B4").Value = 18
$ws.Range("C4").Value = 71
$ws.Range("D4").Value = 28
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 100
$ws.Range("G4").Value = 4935
$ws.Range("H4").Value = 2158
$ws.Range("I4").Value = 384
$ws.Range("J4").Value = 45
$ws.Range("K4").Value = 18
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = "Noun"
